# Auto-generated edit script applying numeric corrections per the Shinryu_Profits diff
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H28").Value = 704.7692
$wsALC.Range("J28").Value = 740.2857
$wsALC.Range("L28").Value = 740.2857
$wsALC.Range("N28").Value = -1710.2857

$wsALC.Range("H43").Value = 649.2
$wsALC.Range("I43").Value = 592
$wsALC.Range("K43").Value = 592
$wsALC.Range("M43").Value = -523

$wsALC.Range("H51").Value = 1716.5
$wsALC.Range("I51").Value = 1999.5
$wsALC.Range("J51").Value = 1575
$wsALC.Range("K51").Value = 1999.5
$wsALC.Range("L51").Value = 1575
$wsALC.Range("M51").Value = -1515.5
$wsALC.Range("N51").Value = -2543

$wsALC.Range("H64").Value = 3977.96
$wsALC.Range("I64").Value = 3833.2917
$wsALC.Range("K64").Value = 3833.2917
$wsALC.Range("M64").Value = -3585.2917

$wsALC.Range("H67").Value = 3977.96
$wsALC.Range("I67").Value = 3833.2917
$wsALC.Range("K67").Value = 3833.2917
$wsALC.Range("M67").Value = -2975.2917

$wsALC.Range("H98").Value = 1938.3334
$wsALC.Range("I98").Value = 1553.4615
$wsALC.Range("J98").Value = 2939
$wsALC.Range("K98").Value = 1553.4615
$wsALC.Range("L98").Value = 2939
$wsALC.Range("M98").Value = -55.46149999999989
$wsALC.Range("N98").Value = -5935

$wsALC.Range("H122").Value = 1938.3334
$wsALC.Range("I122").Value = 1553.4615
$wsALC.Range("J122").Value = 2939
$wsALC.Range("K122").Value = 4660.3845
$wsALC.Range("L122").Value = 8817
$wsALC.Range("M122").Value = -2210.3845
$wsALC.Range("N122").Value = -13717

$wsALC.Range("H125").Value = 894.8889
$wsALC.Range("J125").Value = 0
$wsALC.Range("L125").Value = 0
$wsALC.Range("N125").ClearContents()

$wsALC.Range("H129").Value = 845.54
$wsALC.Range("J129").Value = 845.54
$wsALC.Range("L129").Value = 2536.62
$wsALC.Range("N129").Value = -12536.62

$wsALC.Range("H131").Value = 26070.537
$wsALC.Range("I131").Value = 31207.06
$wsALC.Range("J131").Value = 4882.375
$wsALC.Range("K131").Value = 93621.18000000001
$wsALC.Range("L131").Value = 14647.125
$wsALC.Range("M131").Value = -88581.18000000001
$wsALC.Range("N131").Value = -24727.125

$wsALC.Range("H135").Value = 1209.4445
$wsALC.Range("J135").Value = 2286
$wsALC.Range("L135").Value = 20574
$wsALC.Range("N135").Value = -25644

$wsALC.Range("H137").Value = 61464.94
$wsALC.Range("I137").Value = 3140.6
$wsALC.Range("J137").Value = 85766.75
$wsALC.Range("K137").Value = 9421.799999999999
$wsALC.Range("L137").Value = 257300.25
$wsALC.Range("M137").Value = -6871.799999999999
$wsALC.Range("N137").Value = -262400.25

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 17245340
$wsARM.Range("I32").Value = 20835376
$wsARM.Range("J32").Value = 13163.8
$wsARM.Range("K32").Value = 20835376
$wsARM.Range("L32").Value = 13163.8
$wsARM.Range("M32").Value = -20835089
$wsARM.Range("N32").Value = -13737.8

$wsARM.Range("H61").Value = 5703.6665
$wsARM.Range("I61").Value = 2138.8333
$wsARM.Range("J61").Value = 12833.333
$wsARM.Range("K61").Value = 2138.8333
$wsARM.Range("L61").Value = 12833.333
$wsARM.Range("M61").Value = -1926.8333
$wsARM.Range("N61").Value = -13257.333

$wsARM.Range("H136").Value = 5703.6665
$wsARM.Range("I136").Value = 2138.8333
$wsARM.Range("J136").Value = 12833.333
$wsARM.Range("K136").Value = 6416.499899999999
$wsARM.Range("L136").Value = 38499.999
$wsARM.Range("M136").Value = -3866.499899999999
$wsARM.Range("N136").Value = -43599.999

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H22").Value = 1100.5
$wsCRP.Range("I22").Value = 1117.1666
$wsCRP.Range("J22").Value = 1075.5
$wsCRP.Range("K22").Value = 1117.1666
$wsCRP.Range("L22").Value = 1075.5
$wsCRP.Range("M22").Value = -767.1666
$wsCRP.Range("N22").Value = -1775.5

$wsCRP.Range("H31").Value = 4022.75
$wsCRP.Range("I31").Value = 3057.5454
$wsCRP.Range("J31").Value = 5202.4443
$wsCRP.Range("K31").Value = 3057.5454
$wsCRP.Range("L31").Value = 5202.4443
$wsCRP.Range("M31").Value = -2762.5454
$wsCRP.Range("N31").Value = -5792.4443

$wsCRP.Range("H34").Value = 4022.75
$wsCRP.Range("I34").Value = 3057.5454
$wsCRP.Range("J34").Value = 5202.4443
$wsCRP.Range("K34").Value = 3057.5454
$wsCRP.Range("L34").Value = 5202.4443
$wsCRP.Range("M34").Value = -2855.5454
$wsCRP.Range("N34").Value = -5606.4443

$wsCRP.Range("H132").Value = 2252.8096
$wsCRP.Range("I132").Value = 1350.6875
$wsCRP.Range("K132").Value = 4052.0625
$wsCRP.Range("M132").Value = -1522.0625

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H86").Value = 1588.8889
$wsCUL.Range("J86").Value = 1588.8889
$wsCUL.Range("L86").Value = 4766.6667
$wsCUL.Range("N86").Value = -7138.6667

$wsCUL.Range("H89").Value = 1588.8889
$wsCUL.Range("J89").Value = 1588.8889
$wsCUL.Range("L89").Value = 14300.0001
$wsCUL.Range("N89").Value = -26156.0001

$wsCUL.Range("H97").Value = 538.4545000000001
$wsCUL.Range("J97").Value = 637.6
$wsCUL.Range("L97").Value = 1912.8
$wsCUL.Range("N97").Value = -2904.8

$wsCUL.Range("H98").Value = 249.5
$wsCUL.Range("I98").Value = 242.28572
$wsCUL.Range("J98").Value = 300
$wsCUL.Range("K98").Value = 726.85716
$wsCUL.Range("L98").Value = 900
$wsCUL.Range("M98").Value = 771.14284
$wsCUL.Range("N98").Value = -3896

$wsCUL.Range("H122").Value = 7575.129
$wsCUL.Range("J122").Value = 1224.9166
$wsCUL.Range("L122").Value = 11024.2494
$wsCUL.Range("N122").Value = -15924.2494

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H22").Value = 769977.1
$wsLTW.Range("I22").Value = 1111588.8
$wsLTW.Range("J22").Value = 1351
$wsLTW.Range("K22").Value = 1111588.8
$wsLTW.Range("L22").Value = 1351
$wsLTW.Range("M22").Value = -1111293.8
$wsLTW.Range("N22").Value = -1941

$wsLTW.Range("H27").Value = 769977.1
$wsLTW.Range("I27").Value = 1111588.8
$wsLTW.Range("J27").Value = 1351
$wsLTW.Range("K27").Value = 1111588.8
$wsLTW.Range("L27").Value = 1351
$wsLTW.Range("M27").Value = -1111481.8
$wsLTW.Range("N27").Value = -1565

$wsLTW.Range("H40").Value = 3898.2273
$wsLTW.Range("I40").Value = 3536.7778
$wsLTW.Range("J40").Value = 5524.75
$wsLTW.Range("K40").Value = 3536.7778
$wsLTW.Range("L40").Value = 5524.75
$wsLTW.Range("M40").Value = -3400.7778
$wsLTW.Range("N40").Value = -5796.75

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H126").Value = 1453.1177
$wsWVR.Range("I126").Value = 1554.2
$wsWVR.Range("J126").Value = 695
$wsWVR.Range("K126").Value = 4662.6
$wsWVR.Range("L126").Value = 2085
$wsWVR.Range("M126").Value = -2192.6
$wsWVR.Range("N126").Value = -7025

$wsWVR.Range("H132").Value = 2526.4075
$wsWVR.Range("I132").Value = 2093.8125
$wsWVR.Range("K132").Value = 6281.4375
$wsWVR.Range("M132").Value = -3751.4375
